# make interval function di sisi mahasiswa saat menampilan data absensi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set column C (TANGGAL) values to a single date value (interval) for rows 2-10
$ws.Range("C2:C10").Value = 43748

# Update the active selection to C6
$ws.Range("C6").Select()
